# Adds a new "26/05/2024" daily-change column (column C) to the
# "Nifty_50" worksheet, next to the existing "05/04/2024" column (B).
# This is the data backing the new HTML Report generation feature.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nifty_50")

# Row 1 holds the column header ("26/05/2024"); rows 2-51 hold the
# per-stock price change for that date, in the same row order as the
# existing "Stock Name" (A) / "05/04/2024" (B) columns.
$newColumnValues = @(
    "26/05/2024",
    "−33.75",
    "−34.15",
    "−2.45",
    "+27.30",
    "−7.15",
    "−10.15",
    "−40.20",
    "+2.00",
    "−15.75",
    "−4.80",
    "+35.60",
    "−32.75",
    "−13.10",
    "−18.90",
    "−3.45",
    "−0.60",
    "−3.65",
    "−0.050",
    "−30.70",
    "−41.05",
    "−4.10",
    "+3.55",
    "+7.25",
    "−7.05",
    "−9.80",
    "+2.15",
    "−6.30",
    "−8.10",
    "−10.95",
    "+23.90",
    "−42.50",
    "−2.10",
    "−0.85",
    "+65.15",
    "+2.70",
    "+7.70",
    "−11.05",
    "−0.80",
    "+2.20",
    "+2.85",
    "+15.45",
    "−18.95",
    "+12.70",
    "+3.00",
    "−26.45",
    "−0.45",
    "−9.80",
    "+67.40",
    "+2.95",
    "−3.45"
)

for ($i = 0; $i -lt $newColumnValues.Length; $i++) {
    $row = $i + 1
    $value = $newColumnValues[$i]

    # Values such as "+27.30" look like numbers to Excel's auto-detection
    # and would otherwise be stored as the number 27.3. Prefixing with an
    # apostrophe forces them to stay literal text, matching the existing
    # "+"/"−" prefixed text values already used in column B.
    if ($value.StartsWith("+")) {
        $ws.Cells.Item($row, 3).Value = "'" + $value
    } else {
        $ws.Cells.Item($row, 3).Value = $value
    }
}
